# Hide Slide 18 of the presentation (Session 10 slide deck).
# This mirrors the XML change `show="0"` being added to the <p:sld>
# root element of ppt/slides/slide18.xml, i.e. marking the slide as
# hidden in the slide show.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(18)
$s.SlideShowTransition.Hidden = $true
